{"js": "// Locate the \"Dear {{ users[0]. other_parties_contact }}\" paragraph by its\n// text, append a trailing comma to it, then insert the new \"job info\"\n// paragraph (plus the blank spacer paragraphs that came with it) right\n// after, exactly as described by the commit's diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet dearParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Dear\") !== -1 && t.indexOf(\"other_parties_contact\") !== -1) {\n    dearParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!dearParagraph) {\n  throw new Error(\"Could not find the 'Dear ...' paragraph\");\n}\n\n// Append the trailing comma to the end of the \"Dear\" paragraph.\ndearParagraph.getRange(\"End\").insertText(\",\", \"Before\");\n\n// Blank paragraph right after \"Dear ...,\"\nlet p = dearParagraph.insertParagraph(\"\", \"After\");\n\n// New paragraph introducing the job title / employer / start date.\np = p.insertParagraph(\n  \"I am a {{ users[0].job_title }}with {{ users[0].other_parties }}. I have been in this position since {{ users[0].employment_start_date }}.\",\n  \"After\"\n);\n\n// Three trailing blank paragraphs added along with the new text.\np = p.insertParagraph(\"\", \"After\");\np = p.insertParagraph(\"\", \"After\");\np = p.insertParagraph(\"\", \"After\");\n\nawait context.sync();\n", "ps1": "# Locate the \"Dear {{ users[0]. other_parties_contact }}\" paragraph, append a\n# trailing comma to it, then insert the new \"job info\" paragraph (plus the\n# blank spacer paragraphs that came along with it) right after - matching the\n# commit's diff.\n\n$d = $word.ActiveDocument\n\n$dearIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if ($t -like \"*Dear*\" -and $t -like \"*other_parties_contact*\") {\n        $dearIndex = $i\n        break\n    }\n}\n\nif ($dearIndex -eq -1) {\n    throw \"Could not find the 'Dear ...' paragraph\"\n}\n\n$dearParagraph = $d.Paragraphs($dearIndex)\n$dearRange = $dearParagraph.Range\n$dearRange.MoveEnd(1, -1) | Out-Null\n$dearRange.Collapse(0) | Out-Null\n$dearRange.InsertAfter(\",\")\n\n# Blank paragraph right after \"Dear ...,\"\n$tailRange = $dearParagraph.Range\n$tailRange.Collapse(0) | Out-Null\n$tailRange.InsertParagraphAfter()\n\n# New paragraph introducing the job title / employer / start date.\n$jobParagraph = $d.Paragraphs($dearIndex + 1)\n$jobRange = $jobParagraph.Range\n$jobRange.Collapse(0) | Out-Null\n$jobRange.InsertParagraphAfter()\n\n$textParagraph = $d.Paragraphs($dearIndex + 2)\n$textRange = $textParagraph.Range\n$textRange.MoveEnd(1, -1) | Out-Null\n$textRange.InsertAfter(\"I am a {{ users[0].job_title }}with {{ users[0].other_parties }}. I have been in this position since {{ users[0].employment_start_date }}.\")\n\n# Three trailing blank paragraphs added along with the new text.\n$lastRange = $textParagraph.Range\n$lastRange.Collapse(0) | Out-Null\n$lastRange.InsertParagraphAfter()\n\n$lastParagraph = $d.Paragraphs($dearIndex + 3)\n$lastRange = $lastParagraph.Range\n$lastRange.Collapse(0) | Out-Null\n$lastRange.InsertParagraphAfter()\n\n$lastParagraph = $d.Paragraphs($dearIndex + 4)\n$lastRange = $lastParagraph.Range\n$lastRange.Collapse(0) | Out-Null\n$lastRange.InsertParagraphAfter()\n"}
